$wb = $excel.ActiveWorkbook

# --- 1. Create the new "Player Info" sheet. We copy the "ODI Batting"
#        sheet so the new sheet inherits identical sheetPr/pageMargins/
#        header styling, then trim it down to a 4-column x 2-row table
#        and overwrite its contents. The copy is placed immediately
#        before "ODI Batting", i.e. it becomes the first sheet.
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Copy($batting)

# NOTE: sheet object references drift once the sheet collection changes
# (copy/add/move/rename), so re-fetch by position instead of reusing
# $batting - Copy(Before) always drops the new copy immediately in front
# of the template sheet, i.e. at the very first tab position.
$playerInfo = $wb.Worksheets.Item(1)
$playerInfo.Name = "Player Info"

# Trim the copied sheet down to A1:D2 (drop extra rows/cols), keeping the
# bold/bordered header style that lived on row 1 of the template sheet.
$playerInfo.Range("3:5").Delete()
$playerInfo.Range("E:J").Delete()

# Headers
$playerInfo.Cells.Item(1,1).Value = "ID"
$playerInfo.Cells.Item(1,2).Value = "NAME"
$playerInfo.Cells.Item(1,3).Value = "BATTING_HAND"
$playerInfo.Cells.Item(1,4).Value = "BOWL_STYLE"

# Data row (ID kept as text, like the rest of the workbook's values)
$playerInfo.Cells.Item(2,1).NumberFormat = "@"
$playerInfo.Cells.Item(2,1).Value = "4473"
$playerInfo.Cells.Item(2,1).ClearFormats()
$playerInfo.Cells.Item(2,2).Value = "Binura Fernando"
$playerInfo.Cells.Item(2,3).Value = "Right Handed"
$playerInfo.Cells.Item(2,4).Value = "Left Arm Medium Fast"

# --- 2. ODI Batting: rename MATCH_CARD_LINK -> MATCH_CODE, replace the
#        full scorecard URLs with just the numeric match code.
$batting = $wb.Worksheets.Item("ODI Batting")

$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4465"
$batting.Range("D2").ClearFormats()

$batting.Range("D3").NumberFormat = "@"
$batting.Range("D3").Value = "4469"
$batting.Range("D3").ClearFormats()

$batting.Range("D4").NumberFormat = "@"
$batting.Range("D4").Value = "4470"
$batting.Range("D4").ClearFormats()

$batting.Range("D5").NumberFormat = "@"
$batting.Range("D5").Value = "4471"
$batting.Range("D5").ClearFormats()

# --- 3. ODI Bowling: same rename + value simplification, but column B.
$bowling = $wb.Worksheets.Item("ODI Bowling")

$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4465"
$bowling.Range("B2").ClearFormats()

$bowling.Range("B3").NumberFormat = "@"
$bowling.Range("B3").Value = "4469"
$bowling.Range("B3").ClearFormats()

$bowling.Range("B4").NumberFormat = "@"
$bowling.Range("B4").Value = "4470"
$bowling.Range("B4").ClearFormats()

Write-Host "done"
